$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A2=101, B2=9, C2=30, D2=15, E2=60, F2=15
$ws.Range("A2").Value = 101
$ws.Range("B2").Value = 9
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 60
$ws.Range("F2").Value = 15

# Row 3: A3=501, B3=9, C3=52, D3=30, E3=75, F3=45
$ws.Range("A3").Value = 501
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 52
$ws.Range("D3").Value = 30
$ws.Range("E3").Value = 75
$ws.Range("F3").Value = 45

# Row 4: A4=701, B4=3, C4=90, D4=45, E4=97, F4=15
$ws.Range("A4").Value = 701
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 90
$ws.Range("D4").Value = 45
$ws.Range("E4").Value = 97
$ws.Range("F4").Value = 15

# Row 5: A5=201, B5=9, C5=30, D5=15, E5=45, F5=30
$ws.Range("A5").Value = 201
$ws.Range("B5").Value = 9
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 15
$ws.Range("E5").Value = 45
$ws.Range("F5").Value = 30

# Row 6: A6=902, B6=1, C6=0, D6=0, E6=0, F6=0
$ws.Range("A6").Value = 902
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0

# Row 7: A7=301, B7=6, C7=45, D7=30, E7=60, F7=45
$ws.Range("A7").Value = 301
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 45
$ws.Range("D7").Value = 30
$ws.Range("E7").Value = 60
$ws.Range("F7").Value = 45

# Row 8: A8=401, B8=9, C8=48, D8=67, E8=75, F8=45
$ws.Range("A8").Value = 401
$ws.Range("B8").Value = 9
$ws.Range("C8").Value = 48
$ws.Range("D8").Value = 67
$ws.Range("E8").Value = 75
$ws.Range("F8").Value = 45

# Row 9: A9=601, B9=9, C9=60, D9=67, E9=60, F9=42
$ws.Range("A9").Value = 601
$ws.Range("B9").Value = 9
$ws.Range("C9").Value = 60
$ws.Range("D9").Value = 67
$ws.Range("E9").Value = 60
$ws.Range("F9").Value = 42

# Row 10: A10=1203, B10=3, C10=15, D10=15, E10=15, F10=15
$ws.Range("A10").Value = 1203
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 15
$ws.Range("D10").Value = 15
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 15

# Row 11: A11=901, B11=16, C11=15, D11=45, E11=60, F11=60
$ws.Range("A11").Value = 901
$ws.Range("B11").Value = 16
$ws.Range("C11").Value = 15
$ws.Range("D11").Value = 45
$ws.Range("E11").Value = 60
$ws.Range("F11").Value = 60

# Row 12: A12=1001, B12=18, C12=30, D12=75, E12=60, F12=72
$ws.Range("A12").Value = 1001
$ws.Range("B12").Value = 18
$ws.Range("C12").Value = 30
$ws.Range("D12").Value = 75
$ws.Range("E12").Value = 60
$ws.Range("F12").Value = 72

# Row 13: A13=801, B13=3, C13=67, D13=65, E13=52, F13=45
$ws.Range("A13").Value = 801
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = 67
$ws.Range("D13").Value = 65
$ws.Range("E13").Value = 52
$ws.Range("F13").Value = 45

# Row 14: A14=1201, B14=2, C14=10, D14=10, E14=10, F14=10
$ws.Range("A14").Value = 1201
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 10
$ws.Range("D14").Value = 10
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = 10

# Row 15: A15=1202, B15=2, C15=10, D15=10, E15=10, F15=10
$ws.Range("A15").Value = 1202
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 10
$ws.Range("D15").Value = 10
$ws.Range("E15").Value = 10
$ws.Range("F15").Value = 10

# Row 16: A16=1, B16=0, C16=2, D16=2, E16=2, F16=2
$ws.Range("A16").Value = 1
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 2

# Row 17: A17=2, B17=0, C17=2, D17=2, E17=2, F17=2
$ws.Range("A17").Value = 2
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 2

# Row 18: A18=1101, B18=0, C18=15, D18=30, E18=30, F18=0
$ws.Range("A18").Value = 1101
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 15
$ws.Range("D18").Value = 30
$ws.Range("E18").Value = 30
$ws.Range("F18").Value = 0

# Row 19: A19=3, B19=0, C19=3, D19=3, E19=3, F19=3
$ws.Range("A19").Value = 3
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 3

# Row 20: A20=502, B20=0, C20=4, D20=0, E20=0, F20=0
$ws.Range("A20").Value = 502
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0

# Row 21: A21=802, B21=0, C21=4, D21=5, E21=4, F21=0
$ws.Range("A21").Value = 802
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 0
